$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: values that look like valid Excel numbers must be
# written with a temporary text number format so they stay strings
# (matching the original inlineStr cell type) instead of being parsed
# into numeric values by Excel.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value2 = '27.335.67'
$ws.Range("E2").Value2 = '  +0.07%  '
$ws.Range("D3").Value2 = '1.653.71'
$ws.Range("E3").Value2 = '  -0.41%  '
$ws.Range("E4").Value2 = '  -0.03%  '
Set-TextValue "D5" '218.01'
$ws.Range("E5").Value2 = '  -0.82%  '
Set-TextValue "D6" '0.512'
$ws.Range("E6").Value2 = '  +0.94%  '
$ws.Range("E7").Value2 = '  +0.02%  '
Set-TextValue "D8" '0.256'
$ws.Range("E8").Value2 = '  -0.06%  '
Set-TextValue "D9" '0.0630'
$ws.Range("E9").Value2 = '  +0.69%  '
Set-TextValue "D10" '20.10'
$ws.Range("E10").Value2 = '  +0.19%  '
Set-TextValue "D11" '0.0850'
$ws.Range("E11").Value2 = '  +0.32%  '
$ws.Range("D12").Value2 = '1.887.00'
$ws.Range("E12").Value2 = '  -0.29%  '
$ws.Range("D13").Value2 = '1.631.56'
$ws.Range("E13").Value2 = '  -1.78%  '
Set-TextValue "D14" '4.13'
$ws.Range("E14").Value2 = '  -1.44%  '
Set-TextValue "D15" '0.543'
$ws.Range("E15").Value2 = '  +1.99%  '
Set-TextValue "D16" '67.83'
$ws.Range("E16").Value2 = '  +0.79%  '
$ws.Range("D17").Value2 = '27.349.95'
$ws.Range("E17").Value2 = '  +0.15%  '
$ws.Range("D18").Value2 = '0.0₃0740'
$ws.Range("E18").Value2 = '  +0.65%  '
Set-TextValue "D19" '220.73'
$ws.Range("E19").Value2 = '  -1.14%  '
$ws.Range("E20").Value2 = '  -0.16%  '
Set-TextValue "D21" '6.84'
$ws.Range("E21").Value2 = '  +1.54%  '
Set-TextValue "D22" '2.54'
$ws.Range("E22").Value2 = '  +4.38%  '
Set-TextValue "D23" '4.44'
$ws.Range("E23").Value2 = '  -0.25%  '
Set-TextValue "D24" '9.23'
$ws.Range("E24").Value2 = '  -0.60%  '
Set-TextValue "D25" '147.40'
$ws.Range("E25").Value2 = '  +0.30%  '
Set-TextValue "D26" '7.55'
$ws.Range("E26").Value2 = '  +1.30%  '
$ws.Range("E27").Value2 = '  +0.01%  '
$ws.Range("E28").Value2 = '  -0.93%  '
Set-TextValue "D29" '15.85'
$ws.Range("E29").Value2 = '  -1.41%  '
$ws.Range("E30").Value2 = '  -1.46%  '
$ws.Range("E31").Value2 = '  -1.02%  '
$ws.Range("E32").Value2 = '  -0.90%  '
$ws.Range("E33").Value2 = '  +0.51%  '
$ws.Range("E34").Value2 = '  +0.91%  '
$ws.Range("D35").Value2 = '1.258.90'
$ws.Range("E35").Value2 = '  -0.30%  '
$ws.Range("E36").Value2 = '  +0.06%  '
$ws.Range("E37").Value2 = '  -0.40%  '
Set-TextValue "D38" '0.545'
$ws.Range("E38").Value2 = '  +1.13%  '
Set-TextValue "D39" '0.843'
$ws.Range("E39").Value2 = '  +0.33%  '
$ws.Range("E40").Value2 = '  -0.09%  '
Set-TextValue "D41" '0.809'
$ws.Range("E42").Value2 = '  +5.05%  '
$ws.Range("E43").Value2 = '  +0.98%  '
$ws.Range("D44").Value2 = '1.796.41'
$ws.Range("E44").Value2 = '  -0.30%  '
Set-TextValue "D45" '62.17'
$ws.Range("E45").Value2 = '  +0.38%  '
Set-TextValue "D46" '91.86'
$ws.Range("E46").Value2 = '  -0.27%  '
$ws.Range("E47").Value2 = '  -0.37%  '
$ws.Range("D48").Value2 = '0.0₆0107'
$ws.Range("E48").Value2 = '  +23.34%  '
$ws.Range("E49").Value2 = '  -0.49%  '
Set-TextValue "D50" '7.69'
$ws.Range("E50").Value2 = '  +0.23%  '
Set-TextValue "D51" '0.0974'
$ws.Range("E51").Value2 = '  -0.79%  '
